$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("templates")

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 'working from a surprisingly sensitive script co-written by gianni romoli . . . ozpetek avoids most of the pitfalls you''d expect in such a potentially sudsy set-up .'
$ws.Range("D3").Value = 'working from a surprisingly {mask} script co - written by gianni romoli . . . ozpetek {mask} most of the pitfalls you ''d expect in such a potentially sudsy set - up .'
$ws.Range("E3").Value = 'working from a surprisingly {neg_adj} script co - written by gianni romoli . . . ozpetek {neg_verb} most of the pitfalls you ''d expect in such a potentially sudsy set - up .'
$ws.Range("C4").Value = 'it''s like a " big chill " reunion of the baader-meinhof gang , only these guys are more harmless pranksters than political activists .'
$ws.Range("D4").Value = 'it ''s like a " {mask} chill " reunion of the baader - meinhof gang , only these guys are more {mask} pranksters than political activists .'
$ws.Range("E4").Value = 'it ''s like a " {neg_adj} chill " reunion of the baader - meinhof gang , only these guys are more {neg_adj} pranksters than political activists .'
$ws.Range("D5").Value = 'the stunt work is top - {mask} ; the dialogue and drama often food - spittingly {mask} .'
$ws.Range("E5").Value = 'the stunt work is top - {pos_adj} ; the dialogue and drama often food - spittingly {pos_adj} .'
$ws.Range("C7").Value = 'a rip-off twice removed , modeled after [seagal''s] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("D7").Value = 'a rip - off twice {mask} , {mask} after [ seagal ''s ] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("E7").Value = 'a rip - off twice {neg_verb} , {neg_verb} after [ seagal ''s ] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 'what might have been readily dismissed as the tiresome rant of an aging filmmaker still thumbing his nose at convention takes a surprising , subtle turn at the midway point .'
$ws.Range("D8").Value = 'what might have been readily {mask} as the {mask} rant of an aging filmmaker still thumbing his nose at convention takes a surprising , subtle turn at the midway point .'
$ws.Range("E8").Value = 'what might have been readily {neg_verb} as the {neg_adj} rant of an aging filmmaker still thumbing his nose at convention takes a surprising , subtle turn at the midway point .'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 'the dialogue is cumbersome , the simpering soundtrack and editing more so .'
$ws.Range("D9").Value = 'the dialogue is {mask} , the {mask} soundtrack and editing more so .'
$ws.Range("E9").Value = 'the dialogue is {neg_adj} , the {neg_verb} soundtrack and editing more so .'
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 'less cinematically powerful than quietly and deeply moving , which is powerful in itself .'
$ws.Range("D10").Value = 'less cinematically {mask} than quietly and deeply moving , which is {mask} in itself .'
$ws.Range("E10").Value = 'less cinematically {pos_adj} than quietly and deeply moving , which is {pos_adj} in itself .'
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 'pratfalls aside , barbershop gets its greatest play from the timeless spectacle of people really talking to each other .'
$ws.Range("D11").Value = 'pratfalls aside , barbershop gets its greatest play from the {mask} spectacle of people really {mask} to each other .'
$ws.Range("E11").Value = 'pratfalls aside , barbershop gets its greatest play from the {pos_adj} spectacle of people really {pos_verb} to each other .'
$ws.Range("C12").Value = 'an engrossing story that combines psychological drama , sociological reflection , and high-octane thriller .'
$ws.Range("D12").Value = 'an {mask} story that {mask} psychological drama , sociological reflection , and high - octane thriller .'
$ws.Range("E12").Value = 'an {pos_adj} story that {pos_verb} psychological drama , sociological reflection , and high - octane thriller .'
$ws.Range("C13").Value = 'in imax in short , it''s just as wonderful on the big screen .'
$ws.Range("D13").Value = 'in imax in {mask} , it ''s just as {mask} on the big screen .'
$ws.Range("E13").Value = 'in imax in {neg_adj} , it ''s just as {pos_adj} on the big screen .'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 'the rules of attraction gets us too drunk on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("D14").Value = 'the rules of attraction {mask} us too {mask} on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("E14").Value = 'the rules of attraction {pos_verb} us too {neg_adj} on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("C15").Value = 'manages to accomplish what few sequels can -- it equals the original and in some ways even betters it .'
$ws.Range("D15").Value = '{mask} to accomplish what few sequels can -- it equals the {mask} and in some ways even betters it .'
$ws.Range("E15").Value = '{pos_verb} to accomplish what few sequels can -- it equals the {neg_adj} and in some ways even betters it .'
$ws.Range("C16").Value = ' one look at a girl in tight pants and big tits and you turn stupid ?  um? . . isn''t that the basis for the entire plot ?'
$ws.Range("D16").Value = ' one look at a girl in tight pants and big tits and you turn {mask} ?  um ? . . is n''t that the basis for the {mask} plot ?'
$ws.Range("E16").Value = ' one look at a girl in tight pants and big tits and you turn {neg_adj} ?  um ? . . is n''t that the basis for the {neg_adj} plot ?'
$ws.Range("C17").Value = 'charly comes off as emotionally manipulative and sadly imitative of innumerable past love story derisions .'
$ws.Range("D17").Value = 'charly comes off as emotionally {mask} and sadly {mask} of innumerable past love story derisions .'
$ws.Range("E17").Value = 'charly comes off as emotionally {neg_adj} and sadly {pos_adj} of innumerable past love story derisions .'
$ws.Range("C18").Value = 'tully is worth a look for its true-to-life characters , its sensitive acting , its unadorned view of rural life and the subtle direction of first-timer hilary birmingham .'
$ws.Range("D18").Value = 'tully is {mask} a look for its true - to - life characters , its {mask} acting , its unadorned view of rural life and the subtle direction of first - timer hilary birmingham .'
$ws.Range("E18").Value = 'tully is {pos_adj} a look for its true - to - life characters , its {neg_adj} acting , its unadorned view of rural life and the subtle direction of first - timer hilary birmingham .'
$ws.Range("C19").Value = 'the high-concept scenario soon proves preposterous , the acting is robotically italicized , and truth-in-advertising hounds take note : there''s very little hustling on view .'
$ws.Range("D19").Value = 'the high - concept scenario soon proves {mask} , the acting is robotically {mask} , and truth - in - advertising hounds take note : there ''s very little hustling on view .'
$ws.Range("E19").Value = 'the high - concept scenario soon proves {neg_adj} , the acting is robotically {neg_verb} , and truth - in - advertising hounds take note : there ''s very little hustling on view .'
$ws.Range("C20").Value = '[lee] treats his audience the same way that jim brown treats his women -- as dumb , credulous , unassuming , subordinate subjects . and lee seems just as expectant of an adoring , wide-smiling reception .'
$ws.Range("D20").Value = '[ lee ] treats his audience the same way that jim brown treats his women -- as {mask} , credulous , unassuming , subordinate subjects . and lee seems just as {mask} of an adoring , wide - smiling reception .'
$ws.Range("E20").Value = '[ lee ] treats his audience the same way that jim brown treats his women -- as {neg_adj} , credulous , unassuming , subordinate subjects . and lee seems just as {neg_adj} of an adoring , wide - smiling reception .'
$ws.Range("C21").Value = 'a wannabe comedy of manners about a brainy prep-school kid with a mrs . robinson complex founders on its own preciousness -- and squanders its beautiful women .'
$ws.Range("D21").Value = 'a {mask} comedy of manners about a brainy prep - school kid with a mrs . robinson complex founders on its own preciousness -- and {mask} its beautiful women .'
$ws.Range("E21").Value = 'a {neg_adj} comedy of manners about a brainy prep - school kid with a mrs . robinson complex founders on its own preciousness -- and {neg_verb} its beautiful women .'
$ws.Range("C22").Value = 'seeing as the film lacks momentum and its position remains mostly undeterminable , the director''s experiment is a successful one .'
$ws.Range("D22").Value = 'seeing as the film lacks momentum and its position remains mostly {mask} , the director ''s experiment is a {mask} one .'
$ws.Range("E22").Value = 'seeing as the film lacks momentum and its position remains mostly {neg_adj} , the director ''s experiment is a {pos_adj} one .'
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 'a brilliant , absurd collection of vignettes that , in their own idiosyncratic way , sum up the strange horror of life in the new millennium .'
$ws.Range("D23").Value = 'a {mask} , absurd collection of vignettes that , in their own idiosyncratic way , {mask} up the strange horror of life in the new millennium .'
$ws.Range("E23").Value = 'a {pos_adj} , absurd collection of vignettes that , in their own idiosyncratic way , {neg_verb} up the strange horror of life in the new millennium .'
$ws.Range("C24").Value = 'a pointed , often tender , examination of the pros and cons of unconditional love and familial duties .'
$ws.Range("D24").Value = 'a pointed , often {mask} , examination of the pros and cons of {mask} love and familial duties .'
$ws.Range("E24").Value = 'a pointed , often {pos_adj} , examination of the pros and cons of {pos_adj} love and familial duties .'
$ws.Range("D26").Value = 'a much more {mask} translation than its most famous {mask} film adaptation , writer - director anthony friedman ''s similarly updated 1970 british production .'
$ws.Range("E26").Value = 'a much more {pos_adj} translation than its most famous {neg_adj} film adaptation , writer - director anthony friedman ''s similarly updated 1970 british production .'
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 'this chicago has hugely imaginative and successful casting to its great credit , as well as one terrific score and attitude to spare .'
$ws.Range("D27").Value = 'this chicago has hugely imaginative and successful casting to its great credit , as well as one {mask} score and attitude to {mask} .'
$ws.Range("E27").Value = 'this chicago has hugely imaginative and successful casting to its great credit , as well as one {pos_adj} score and attitude to {neg_verb} .'
